$wb = $excel.ActiveWorkbook

# --- Add the two new (empty) sheets ------------------------------------
# Target sheet order/ids (from the workbook.xml diff):
#   teachers  sheetId=1 (existing "staffDetails_L2", renamed)
#   subjects  sheetId=4 (new, empty)
#   rooms     sheetId=3 (new, empty)
#   incomes   sheetId=2 (existing "Incomes", renamed)
#
# Excel assigns sheetId sequentially as sheets are created, so adding
# "rooms" before "subjects" makes rooms -> 3 and subjects -> 4, matching
# the target file.

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$roomsSheet = $wb.Worksheets.Add($null, $last)
$roomsSheet.Name = "rooms"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$subjectsSheet = $wb.Worksheets.Add($null, $last)
$subjectsSheet.Name = "subjects"

# --- Reorder sheets: teachers, subjects, rooms, incomes -----------------
$incomesSheet = $wb.Worksheets.Item("Incomes")
$roomsSheet = $wb.Worksheets.Item("rooms")
$roomsSheet.Move($incomesSheet)

$roomsSheet = $wb.Worksheets.Item("rooms")
$subjectsSheet = $wb.Worksheets.Item("subjects")
$subjectsSheet.Move($roomsSheet)

# --- Lowercase all table/tab titles -------------------------------------
$wb.Worksheets.Item("staffDetails_L2").Name = "teachers"
$wb.Worksheets.Item("Incomes").Name = "incomes"
# "subjects" and "rooms" are already lowercase.

# --- Make "teachers" the active/selected sheet (was "incomes") ----------
$teachersSheet = $wb.Worksheets.Item("teachers")
$teachersSheet.Activate()

# --- Restore the view/selection on the teachers sheet --------------------
# Scroll so column B is the left-most visible column, and the frozen pane
# (still split after row 1) scrolls to show row 6 onward, while the
# existing selected cell (P6) is left untouched.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 2
$aw.ScrollRow = 6
$teachersSheet.Range("P6").Select()
